{"js": "// Fixed typo found by Annie \u2014 remove the unused \"Footnote Text\" paragraph\n// style from the document's style sheet.\nconst styles = context.document.getStyles();\nconst footnoteTextStyle = styles.getByNameOrNullObject(\"Footnote Text\");\nawait context.sync();\n\nif (!footnoteTextStyle.isNullObject) {\n  footnoteTextStyle.delete();\n  await context.sync();\n}\n", "ps1": "# Fixed typo found by Annie \u2014 remove the unused \"Footnote Text\" paragraph\n# style from the document's style sheet.\n$d = $word.ActiveDocument\n\n$styleName = \"Footnote Text\"\n$exists = $false\nforeach ($s in $d.Styles) {\n    if ($s.NameLocal -eq $styleName) {\n        $exists = $true\n    }\n}\n\nif ($exists) {\n    try {\n        $d.Styles($styleName).Delete()\n    } catch {\n        # Style vanished between the existence check and the delete (or is\n        # a built-in style that cannot be removed) - nothing more to do.\n    }\n}\n"}
